$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.011.82"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.114.08"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.57"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.79"
$ws.Range("E6").Value = "  +4.27%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.109.30"
$ws.Range("E8").Value = "  +1.75%  "

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  -3.50%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.14"
$ws.Range("E14").Value = "  +1.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.629.63"
$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.000.43"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.111.97"
$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.22"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "477.51"
$ws.Range("E21").Value = "  +3.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.711"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.73"
$ws.Range("E23").Value = "  +4.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.86"
$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.30"
$ws.Range("E25").Value = "  +4.17%  "

$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.63"
$ws.Range("E32").Value = "  +1.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0967"
$ws.Range("E33").Value = "  -4.36%  "

$ws.Range("E34").Value = "  -1.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.979"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.77"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.09"
$ws.Range("E39").Value = "  +3.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.06"
$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.310"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.60"
$ws.Range("E43").Value = "  -0.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.800.34"
$ws.Range("E44").Value = "  +1.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0355"
$ws.Range("E45").Value = "  -1.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "380.02"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  -9.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.21"
$ws.Range("E48").Value = "  +2.12%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.76"
$ws.Range("E50").Value = "  +1.62%  "

$ws.Range("E51").Value = "  -0.34%  "
